# Thailand Premier League workbook update (28-05-2024 07:50)
# - Swap team-name labels for "Sukhothai FC" and "Buriram United" in the
#   shared-string table (their display names were mixed up).
# - Re-sync the affected match rows so the correct match data travels
#   with the correct fixture:
#     * rows 15 <-> 16 (match ids 13/14)
#     * rows 117 <-> 118 (match ids 115/116)
#     * rows 233 -> 235 -> 236 -> 237 -> 238 -> 233 (match ids 231,233,234,235,236)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = "$firstCol$rowA`:$lastCol$rowA"
    $rangeB = "$firstCol$rowB`:$lastCol$rowB"
    $valsA = $sheet.Range($rangeA).Value2
    $valsB = $sheet.Range($rangeB).Value2
    $sheet.Range($rangeA).Value = $valsB
    $sheet.Range($rangeB).Value = $valsA
}

function Rotate-Rows($sheet, [int[]]$rows, $firstCol, $lastCol) {
    # content of $rows[0] moves to $rows[1], $rows[1] -> $rows[2], ...,
    # last row's original content moves back to $rows[0].
    $snapshots = @()
    foreach ($r in $rows) {
        $rng = "$firstCol$r`:$lastCol$r"
        $snapshots += ,$sheet.Range($rng).Value2
    }
    $count = $rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $rows[($i + 1) % $count]
        $rng = "$firstCol$destRow`:$lastCol$destRow"
        $sheet.Range($rng).Value = $snapshots[$i]
    }
}

# Row id (col A) and date (col D) are untouched; only B:AD need to move.
Swap-Rows $ws 15 16 "B" "AD"
Swap-Rows $ws 117 118 "B" "AD"
Rotate-Rows $ws @(233, 235, 236, 237, 238) "B" "AD"
